$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must remain text
# (matches the original inline-string cell type and preserves formatting
# such as trailing zeros, e.g. "1.00" / "28.00" / "9.00")
$numericLookingText = [ordered]@{
    'D5' = '570.41'
    'D6' = '142.15'
    'D10' = '7.55'
    'D15' = '0.125'
    'D20' = '13.89'
    'D21' = '8.99'
    'D22' = '385.28'
    'D24' = '74.11'
    'D29' = '1.00'
    'D30' = '7.39'
    'D31' = '7.96'
    'D35' = '23.52'
    'D37' = '167.62'
    'D40' = '1.48'
    'D41' = '28.00'
    'D44' = '0.779'
    'D45' = '42.16'
    'D50' = '6.82'
    'D51' = '23.06'
}

foreach ($cellRef in $numericLookingText.Keys) {
    $ws.Range($cellRef).NumberFormat = '@'
}
foreach ($cellRef in $numericLookingText.Keys) {
    $ws.Range($cellRef).Value = $numericLookingText[$cellRef]
}

# Remaining cells: plain text / percentage strings, safe to assign directly
$plainText = [ordered]@{
    'D2' = '60.945.11'
    'E2' = '  -0.74%  '
    'D3' = '3.393.04'
    'E3' = '  -1.22%  '
    'E4' = '  -0.02%  '
    'E5' = '  -0.69%  '
    'E6' = '  -2.05%  '
    'D7' = '3.394.21'
    'E7' = '  -1.17%  '
    'E8' = '  +0.05%  '
    'E9' = '  -0.54%  '
    'E10' = '  -1.63%  '
    'E11' = '  -1.53%  '
    'D13' = '3.972.31'
    'E13' = '  -1.22%  '
    'E14' = '  +1.21%  '
    'E15' = '  +1.73%  '
    'E16' = '  -1.09%  '
    'D17' = '3.393.51'
    'E17' = '  -1.31%  '
    'D18' = '61.010.08'
    'E18' = '  -0.74%  '
    'E19' = '  -1.37%  '
    'E20' = '  -2.14%  '
    'E21' = '  -4.33%  '
    'E22' = '  -2.52%  '
    'E23' = '  -1.37%  '
    'E24' = '  +1.25%  '
    'E25' = '  +0.35%  '
    'E26' = '  -4.47%  '
    'D27' = '3.533.03'
    'E27' = '  -1.10%  '
    'E28' = '  -0.11%  '
    'E30' = '  -2.71%  '
    'E31' = '  -2.68%  '
    'E33' = '  -2.87%  '
    'E35' = '  -1.93%  '
    'E36' = '  -0.35%  '
    'E37' = '  +0.04%  '
    'D38' = '3.423.27'
    'E38' = '  -1.10%  '
    'E39' = '  -2.29%  '
    'E40' = '  -4.71%  '
    'E41' = '  +4.46%  '
    'E42' = '  -1.23%  '
    'E43' = '  -0.01%  '
    'E44' = '  -2.41%  '
    'E45' = '  +0.08%  '
    'E46' = '  -1.20%  '
    'E47' = '  -3.31%  '
    'E48' = '  -2.18%  '
    'D49' = '2.487.43'
    'E49' = '  -3.76%  '
    'E50' = '  -1.43%  '
    'E51' = '  +0.01%  '
}

foreach ($cellRef in $plainText.Keys) {
    $ws.Range($cellRef).Value = $plainText[$cellRef]
}
